$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 36: date label "03-11-2021" (as text, matching the style of the
# existing date cells in column A), a 10000 cupo, and a 0 in column D.
$cellA = $ws.Cells.Item(36, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "03-11-2021"
$cellA.Style = "Normal"

$ws.Cells.Item(36, 2).Value = 10000
$ws.Cells.Item(36, 4).Value = 0
